$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.915.78"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.414.63"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "253.91"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "662.85"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "1.49"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "1.04"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "3.411.89"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("D13").Value = "44.25"
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").Value = "97.693.92"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "6.14"
$ws.Range("E15").Value = "  -4.58%  "
$ws.Range("D16").Value = "'0.0000259"
$ws.Range("E16").Value = "  -3.62%  "
$ws.Range("D17").Value = "4.043.59"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "9.26"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "3.425.85"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "18.31"
$ws.Range("E20").Value = "  +2.81%  "
$ws.Range("D21").Value = "0.526"
$ws.Range("E21").Value = "  -10.95%  "
$ws.Range("D22").Value = "'11.50"
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("D23").Value = "512.45"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "3.44"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").Value = "'0.0000202"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("D27").Value = "'97.10"
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").Value = "'12.50"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").Value = "3.562.93"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "11.82"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "0.144"
$ws.Range("E31").Value = "  -5.52%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("D34").Value = "2.67"
$ws.Range("E34").Value = "  +6.60%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "0.567"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "29.39"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").Value = "7.98"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").Value = "1.49"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").Value = "527.49"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").Value = "0.153"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.871"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").Value = "24.39"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "0.0428"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").Value = "3.68"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "5.66"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "8.67"
$ws.Range("E49").Value = "  -3.97%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "2.23"
$ws.Range("E50").Value = "  +4.94%  "
$ws.Range("D51").Value = "55.69"
$ws.Range("E51").Value = "  +3.00%  "
